# Separate dividends and tax events: rename the "Withheld Tax on Dividends (...)"
# labels on the "Tax Withholding" sheet to "Tax Withholding (...)", and update the
# USD balances on the "Foreign Currencies" sheet that shift as a consequence of
# splitting dividend/tax events during the Schwab export & event generation.

$wb = $excel.ActiveWorkbook

# --- Tax Withholding sheet: shorten/rename the comment labels ---
$wsTax = $wb.Worksheets.Item("Tax Withholding")
$wsTax.Range("B2").Value = "Tax Withholding (NVDA)"
$wsTax.Range("B3").Value = "Tax Withholding (APPL)"

# Column B shrinks now that the text is shorter - re-fit it to the new content.
$wsTax.Columns.Item(2).ColumnWidth = 21.67

# --- Foreign Currencies sheet: updated USD amounts ---
$wsFx = $wb.Worksheets.Item("Foreign Currencies")
$wsFx.Range("B2").Value = 1217.91
$wsFx.Range("B3").Value = 100
$wsFx.Range("B4").Value = 100
